# Updates cryptos list price/volume figures (GitHub Actions scrape refresh).
# Price cells (column D) sometimes look like plain numbers (e.g. "1.00",
# "19.60"); Excel would silently coerce those to numeric values and drop
# the formatting digits (1.00 -> 1), so for any D-column value that parses
# as a number we briefly force a Text number format, assign the literal
# string, then reset the style back to "Normal" so no stray style/format
# is left behind on the cell (matches the original workbook's styling).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.613.17"
$ws.Range("E2").Value = "  +6.30%  "
$ws.Range("D3").Value = "1.941.95"
$ws.Range("E3").Value = "  +2.85%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "251.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.693"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.42%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "48.52"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +12.46%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.381"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +7.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "58.78"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.19%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0773"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.81%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.77"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +13.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.846"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +9.30%  "
$ws.Range("D15").Value = "2.219.15"
$ws.Range("E15").Value = "  +2.56%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.20"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.71%  "
$ws.Range("D17").Value = "1.936.12"
$ws.Range("E17").Value = "  +2.49%  "
$ws.Range("D18").Value = "37.628.37"
$ws.Range("E18").Value = "  +6.27%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "75.64"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.97%  "
$ws.Range("D20").Value = "0.0₃0865"
$ws.Range("E20").Value = "  +4.72%  "
$ws.Range("E21").Value = "  +7.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "254.42"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.82%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.25"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.77%  "
$ws.Range("E24").Value = "  +0.15%  "
$ws.Range("E25").Value = "  -6.34%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "169.37"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.57%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.14"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.94%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.99"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.76%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.94"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.53%  "
$ws.Range("E30").Value = "  +1.26%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.62"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.72%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0618"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.81%  "
$ws.Range("E33").Value = "  +27.44%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.36"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.11%  "
$ws.Range("B35").Value = "Gas"
$ws.Range("C35").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "19.60"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +42.35%  "
$ws.Range("B36").Value = "BinanceUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.08%  "
$ws.Range("B37").Value = "WEMIXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.88"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.64%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.900"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.20%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.46"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.89%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.99%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "106.02"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.85%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0229"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.61"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.18%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.87"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +19.72%  "
$ws.Range("E45").Value = "  +3.89%  "
$ws.Range("D46").Value = "1.360.59"
$ws.Range("E46").Value = "  +2.57%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.44"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.91%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0846"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.38%  "
$ws.Range("E49").Value = "  +2.58%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.86"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +15.99%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.47"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.71%  "
